$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B50/C50 (EOS -> BabyDogeCoin) and B51/C51 (ThetaToken -> EOS)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"

# Update D/E (price, volume) columns for all data rows 2-51
$ws.Range("D2").Value = "'30.379.89"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "'2.109.17"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").Value = "'344.83"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "'0.5233"
$ws.Range("E7").Value = "  +1.39%  "
$ws.Range("D8").Value = "'0.4450"
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "'54.24"
$ws.Range("E9").Value = "  +4.36%  "
$ws.Range("D10").Value = "'0.09384"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").Value = "'1.175"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'25.30"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'8.690"
$ws.Range("E13").Value = "  +6.24%  "
$ws.Range("D14").Value = "'2.119.00"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "'6.927"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "'101.85"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "'0.00001163"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "'1.005"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'21.35"
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").Value = "'0.06720"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Value = "'6.302"
$ws.Range("E21").Value = "  +1.94%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'30.425.18"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").Value = "'12.65"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'2.366.12"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").Value = "'22.04"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").Value = "'2.545"
$ws.Range("E28").Value = "  +1.04%  "
$ws.Range("D29").Value = "'162.17"
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("D30").Value = "'133.44"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'1.147"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "'1.761"
$ws.Range("E32").Value = "  +8.04%  "
$ws.Range("D33").Value = "'0.1057"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "'6.779"
$ws.Range("E34").Value = "  +12.41%  "
$ws.Range("D35").Value = "'6.263"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").Value = "'3.937"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'10.58"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("D38").Value = "'0.02641"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").Value = "'0.06834"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").Value = "'0.7067"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").Value = "'12.59"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").Value = "'1.342"
$ws.Range("E42").Value = "  +3.88%  "
$ws.Range("D43").Value = "'0.2235"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").Value = "'0.6846"
$ws.Range("E44").Value = "  +3.72%  "
$ws.Range("D45").Value = "'14.56"
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("D46").Value = "'2.371"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "'1.392"
$ws.Range("E48").Value = "  +19.43%  "
$ws.Range("D49").Value = "'3.647"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("D50").Value = "'0.00000000344"
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").Value = "'1.224"
$ws.Range("E51").Value = "  +0.62%  "
